$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.24434529525963
$ws.Range("C2").Value = 15.98798310421177
$ws.Range("D2").Value = 3.901488561645881
$ws.Range("E2").Value = 29.33630424589722
$ws.Range("F2").Value = 21.19244645375618
$ws.Range("G2").Value = 2.068328951829761
$ws.Range("H2").Value = 3.23582982603127
$ws.Range("I2").Value = 3.346455671044311
$ws.Range("P2").Value = 11.7678710571247
$ws.Range("Q2").Value = 16.6519096887458
$ws.Range("B3").Value = 19.83404491452785
$ws.Range("C3").Value = 15.06083541376267
$ws.Range("D3").Value = 3.796144748385419
$ws.Range("E3").Value = 27.45330212399358
$ws.Range("F3").Value = 20.3309374745291
$ws.Range("G3").Value = 2.07297968463036
$ws.Range("H3").Value = 2.982130108109396
$ws.Range("I3").Value = 3.168475253322779
$ws.Range("P3").Value = 11.85923295924826
$ws.Range("Q3").Value = 16.1666207780694
$ws.Range("B4").Value = 18.91300103090603
$ws.Range("C4").Value = 14.46560476936423
$ws.Range("D4").Value = 3.729274292747749
$ws.Range("E4").Value = 26.23278851326196
$ws.Range("F4").Value = 19.792866582613
$ws.Range("G4").Value = 2.07592495572274
$ws.Range("H4").Value = 2.820997640627038
$ws.Range("I4").Value = 3.056221160496542
$ws.Range("P4").Value = 11.91791939148089
$ws.Range("Q4").Value = 15.86810596097616
$ws.Range("B5").Value = 18.51737395584194
$ws.Range("C5").Value = 14.22662462568218
$ws.Range("D5").Value = 3.702647104088578
$ws.Range("E5").Value = 25.71893596989928
$ws.Range("F5").Value = 19.56228696292059
$ws.Range("G5").Value = 2.077155587032032
$ws.Range("H5").Value = 2.75378960930067
$ws.Range("I5").Value = 3.010004753504221
$ws.Range("P5").Value = 11.943404496475
$ws.Range("Q5").Value = 15.73893713464746
$ws.Range("B6").Value = 18.44304980382731
$ws.Range("C6").Value = 14.198826641894
$ws.Range("D6").Value = 3.699614661890565
$ws.Range("E6").Value = 25.63228342230575
$ws.Range("F6").Value = 19.51249521288196
$ws.Range("G6").Value = 2.077370327591373
$ws.Range("H6").Value = 2.742335900622129
$ws.Range("I6").Value = 3.002823975041709
$ws.Range("P6").Value = 11.94882197734478
$ws.Range("Q6").Value = 15.70830475881949
$ws.Range("B7").Value = 18.88696306560646
$ws.Range("C7").Value = 14.49524762205506
$ws.Range("D7").Value = 3.732781708952662
$ws.Range("E7").Value = 26.2249901650297
$ws.Range("F7").Value = 19.75907531799321
$ws.Range("G7").Value = 2.075965619517161
$ws.Range("H7").Value = 2.819531202845754
$ws.Range("I7").Value = 3.056754298870644
$ws.Range("P7").Value = 11.92148074146381
$ws.Range("Q7").Value = 15.841509031171
$ws.Range("B8").Value = 20.74412174701266
$ws.Range("C8").Value = 15.71431410343474
$ws.Range("D8").Value = 3.870664877140791
$ws.Range("E8").Value = 28.69922864952596
$ws.Range("F8").Value = 20.85879796876124
$ws.Range("G8").Value = 2.069944485621055
$ws.Range("H8").Value = 3.14873931336443
$ws.Range("I8").Value = 3.286764084671139
$ws.Range("P8").Value = 11.80347848797396
$ws.Range("Q8").Value = 16.45301393072154
$ws.Range("B9").Value = 23.9847716862231
$ws.Range("C9").Value = 17.84333741212275
$ws.Range("D9").Value = 4.116343728166592
$ws.Range("E9").Value = 33.04398820311501
$ws.Range("F9").Value = 22.98345331535855
$ws.Range("G9").Value = 2.058785775659262
$ws.Range("H9").Value = 3.753920340167167
$ws.Range("I9").Value = 3.714043418242903
$ws.Range("P9").Value = 11.58678339590977
$ws.Range("Q9").Value = 17.68518224625051
$ws.Range("B10").Value = 26.03871136407511
$ws.Range("C10").Value = 19.23416383004324
$ws.Range("D10").Value = 4.311185115942086
$ws.Range("E10").Value = 35.05051666660535
$ws.Range("F10").Value = 24.32600997702151
$ws.Range("G10").Value = 2.051226890993335
$ws.Range("H10").Value = 4.137073897476883
$ws.Range("I10").Value = 4.007321852824707
$ws.Range("P10").Value = 11.46189746524951
$ws.Range("Q10").Value = 18.44890189708534
$ws.Range("B11").Value = 26.39814353220931
$ws.Range("C11").Value = 19.39037959456375
$ws.Range("D11").Value = 4.617656838781644
$ws.Range("E11").Value = 28.4413378963673
$ws.Range("F11").Value = 23.65465203734853
$ws.Range("G11").Value = 2.0498429246093
$ws.Range("H11").Value = 4.492300152688726
$ws.Range("I11").Value = 4.071110765071338
$ws.Range("P11").Value = 11.58943319410997
$ws.Range("Q11").Value = 17.73242486140492
$ws.Range("B12").Value = 26.31017654352046
$ws.Range("C12").Value = 19.20337017689211
$ws.Range("D12").Value = 4.838017044454928
$ws.Range("E12").Value = 22.48597269447575
$ws.Range("F12").Value = 22.83200781241673
$ws.Range("G12").Value = 2.050003634790813
$ws.Range("H12").Value = 5.297414652992155
$ws.Range("I12").Value = 4.068682360689743
$ws.Range("P12").Value = 11.72926721177756
$ws.Range("Q12").Value = 16.98624232278696
$ws.Range("B13").Value = 25.84726470326447
$ws.Range("C13").Value = 18.77667581959529
$ws.Range("D13").Value = 5.01237998862362
$ws.Range("E13").Value = 16.56201812557444
$ws.Range("F13").Value = 21.78085996614549
$ws.Range("G13").Value = 2.051387570047776
$ws.Range("H13").Value = 6.317432706305129
$ws.Range("I13").Value = 4.017252459079821
$ws.Range("P13").Value = 11.88680967195127
$ws.Range("Q13").Value = 16.12346427244232
$ws.Range("B14").Value = 25.34308545040505
$ws.Range("C14").Value = 18.37350547551109
$ws.Range("D14").Value = 5.116182715614461
$ws.Range("E14").Value = 12.49264091669018
$ws.Range("F14").Value = 20.92278044868427
$ws.Range("G14").Value = 2.052888961639587
$ws.Range("H14").Value = 7.116731744846066
$ws.Range("I14").Value = 3.960504999053573
$ws.Range("P14").Value = 12.00576682982438
$ws.Range("Q14").Value = 15.45461220619254
$ws.Range("B15").Value = 25.13125716500807
$ws.Range("C15").Value = 18.23190355643588
$ws.Range("D15").Value = 5.132674718159113
$ws.Range("E15").Value = 11.51593739001252
$ws.Range("F15").Value = 20.65438187367861
$ws.Range("G15").Value = 2.053551800275636
$ws.Range("H15").Value = 7.300747388952847
$ws.Range("I15").Value = 3.9361343313207
$ws.Range("P15").Value = 12.0386306197751
$ws.Range("Q15").Value = 15.25954650730678
$ws.Range("B16").Value = 24.34195033592232
$ws.Range("C16").Value = 17.73640174565618
$ws.Range("D16").Value = 5.022720057685066
$ws.Range("E16").Value = 11.39309845988134
$ws.Range("F16").Value = 20.22206244717581
$ws.Range("G16").Value = 2.056465447626719
$ws.Range("H16").Value = 7.02225748233182
$ws.Range("I16").Value = 3.824142457375225
$ws.Range("P16").Value = 12.04649673184934
$ws.Range("Q16").Value = 15.06049289056866
$ws.Range("B17").Value = 24.00371943281733
$ws.Range("C17").Value = 17.57351854959033
$ws.Range("D17").Value = 4.880251754739056
$ws.Range("E17").Value = 13.4979010985307
$ws.Range("F17").Value = 20.36296271496995
$ws.Range("G17").Value = 2.057930077344666
$ws.Range("H17").Value = 6.332660876888924
$ws.Range("I17").Value = 3.767659763510585
$ws.Range("P17").Value = 11.99182856646824
$ws.Range("Q17").Value = 15.27309957087339
$ws.Range("B18").Value = 24.044348790412
$ws.Range("C18").Value = 17.66693945972051
$ws.Range("D18").Value = 4.693155499359199
$ws.Range("E18").Value = 18.0617373087299
$ws.Range("F18").Value = 21.03231258682181
$ws.Range("G18").Value = 2.058178860574236
$ws.Range("H18").Value = 5.299013535784304
$ws.Range("I18").Value = 3.754025291857811
$ws.Range("P18").Value = 11.87671840079358
$ws.Range("Q18").Value = 15.87749421297165
$ws.Range("B19").Value = 24.35879215894935
$ws.Range("C19").Value = 18.01543152623328
$ws.Range("D19").Value = 4.497482009348721
$ws.Range("E19").Value = 24.34212166184587
$ws.Range("F19").Value = 22.01514394023083
$ws.Range("G19").Value = 2.057304220750735
$ws.Range("H19").Value = 4.323344623511946
$ws.Range("I19").Value = 3.785884445289801
$ws.Range("P19").Value = 11.73771417547752
$ws.Range("Q19").Value = 16.70409758324919
$ws.Range("B20").Value = 25.46361636489998
$ws.Range("C20").Value = 18.95625835189501
$ws.Range("D20").Value = 4.273331612968068
$ws.Range("E20").Value = 34.4935351717964
$ws.Range("F20").Value = 23.8909754692182
$ws.Range("G20").Value = 2.053260626197614
$ws.Range("H20").Value = 4.03386592142914
$ws.Range("I20").Value = 3.935342292896853
$ws.Range("P20").Value = 11.50748466950926
$ws.Range("Q20").Value = 18.17589669579615
$ws.Range("B21").Value = 27.06302406206178
$ws.Range("C21").Value = 20.05080150246729
$ws.Range("D21").Value = 4.378876214376258
$ws.Range("E21").Value = 37.30861685159484
$ws.Range("F21").Value = 25.14104151860499
$ws.Range("G21").Value = 2.047180558833289
$ws.Range("H21").Value = 4.380770126744784
$ws.Range("I21").Value = 4.169007924582534
$ws.Range("P21").Value = 11.38581223215721
$ws.Range("Q21").Value = 18.95822593577839
$ws.Range("B22").Value = 28.05083845722103
$ws.Range("C22").Value = 20.6820002794959
$ws.Range("D22").Value = 4.453456234324586
$ws.Range("E22").Value = 38.6401310496351
$ws.Range("F22").Value = 25.90144730192458
$ws.Range("G22").Value = 2.043365079700162
$ws.Range("H22").Value = 4.584481656539372
$ws.Range("I22").Value = 4.313709910463947
$ws.Range("P22").Value = 11.31264952474652
$ws.Range("Q22").Value = 19.43138575637101
$ws.Range("B23").Value = 27.54557798507469
$ws.Range("C23").Value = 20.3181632356919
$ws.Range("D23").Value = 4.409673853925127
$ws.Range("E23").Value = 37.93528199551817
$ws.Range("F23").Value = 25.52562432809528
$ws.Range("G23").Value = 2.045373237940637
$ws.Range("H23").Value = 4.476635571788233
$ws.Range("I23").Value = 4.234884081058053
$ws.Range("P23").Value = 11.34622501537988
$ws.Range("Q23").Value = 19.20377684041285
$ws.Range("B24").Value = 25.51262758496888
$ws.Range("C24").Value = 18.93701595186494
$ws.Range("D24").Value = 4.245717152561948
$ws.Range("E24").Value = 35.15552401392027
$ws.Range("F24").Value = 24.03099774490324
$ws.Range("G24").Value = 2.053131939213595
$ws.Range("H24").Value = 4.059945008673057
$ws.Range("I24").Value = 3.934995918583402
$ws.Range("P24").Value = 11.48737804446741
$ws.Range("Q24").Value = 18.29734762139054
$ws.Range("B25").Value = 23.12287659847034
$ws.Range("C25").Value = 17.3406112688562
$ws.Range("D25").Value = 4.058593364205172
$ws.Range("E25").Value = 31.92129429900232
$ws.Range("F25").Value = 22.36976730959418
$ws.Range("G25").Value = 2.061772739636632
$ws.Range("H25").Value = 3.593540973059603
$ws.Range("I25").Value = 3.602901949824975
$ws.Range("P25").Value = 11.65009461266547
$ws.Range("Q25").Value = 17.31236120644557
